# "9th Stab - Cosmetic Changes"
#
# The report rotates its rolling history columns: the two most-recent
# "UN" rating columns (currently B = Jun_13, C = Jun_10) get pushed out
# to the right (D = Jun_13, E = Jun_10), and two brand-new "UN" columns
# are introduced in their place (B = Jun_17, C = Jun_15) filled with the
# default "UN" rating value for every data row.
#
# Any individual non-default cell (e.g. a highlighted upgrade/downgrade
# note) travels together with the column it originally belonged to.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlNone = -4142
$lastRow = 27

# --- 1. Shift the existing data right: old col B -> D, old col C -> E ---
for ($r = 1; $r -le $lastRow; $r++) {
    $srcB = $ws.Cells.Item($r, 2)
    $srcC = $ws.Cells.Item($r, 3)
    $dstD = $ws.Cells.Item($r, 4)
    $dstE = $ws.Cells.Item($r, 5)

    $dstD.Value = $srcB.Text
    if ($srcB.Interior.Pattern -ne $xlNone) {
        $dstD.Interior.Color = $srcB.Interior.Color
    }

    $dstE.Value = $srcC.Text
    if ($srcC.Interior.Pattern -ne $xlNone) {
        $dstE.Interior.Color = $srcC.Interior.Color
    }
}

# --- 2. Populate the two new leading columns with plain default formatting ---
$ws.Cells.Item(1, 2).ClearFormats()
$ws.Cells.Item(1, 2).Value = "Jun_17"
$ws.Cells.Item(1, 3).ClearFormats()
$ws.Cells.Item(1, 3).Value = "Jun_15"

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).ClearFormats()
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).ClearFormats()
    $ws.Cells.Item($r, 3).Value = "UN"
}

# --- 3. Match the column widths used by the rest of the history columns ---
$ws.Columns.Item(4).ColumnWidth = 45.0
$ws.Columns.Item(5).ColumnWidth = 45.0
